$wb = $excel.ActiveWorkbook

# ----- Worksheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 21
$ws.Range("H21").Value = 20196.2
$ws.Range("J21").Value = 18333.334
$ws.Range("L21").Value = 18333.334
$ws.Range("N21").Value = -19269.334

# Row 23
$ws.Range("H23").Value = 20196.2
$ws.Range("J23").Value = 18333.334
$ws.Range("L23").Value = 18333.334
$ws.Range("N23").Value = -18801.334

# Row 29
$ws.Range("H29").Value = 493.2857
$ws.Range("I29").Value = 530.6
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 1591.8
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = -1310.8
$ws.Range("N29").Value = -1762

# Row 69
$ws.Range("H69").Value = 3794.2
$ws.Range("I69").Value = 3012.9167
$ws.Range("J69").Value = 4515.385
$ws.Range("K69").Value = 9038.750100000001
$ws.Range("L69").Value = 13546.155
$ws.Range("M69").Value = -8164.750100000001
$ws.Range("N69").Value = -15294.155

# Row 70
$ws.Range("H70").Value = 3608.8
$ws.Range("I70").Value = 2117.6
$ws.Range("J70").Value = 5597.067
$ws.Range("K70").Value = 6352.799999999999
$ws.Range("L70").Value = 16791.201
$ws.Range("M70").Value = -6082.799999999999
$ws.Range("N70").Value = -17331.201

# Row 72
$ws.Range("H72").Value = 3794.2
$ws.Range("I72").Value = 3012.9167
$ws.Range("J72").Value = 4515.385
$ws.Range("K72").Value = 27116.2503
$ws.Range("L72").Value = 40638.465
$ws.Range("M72").Value = -22748.2503
$ws.Range("N72").Value = -49374.465

# Row 73
$ws.Range("H73").Value = 3608.8
$ws.Range("I73").Value = 2117.6
$ws.Range("J73").Value = 5597.067
$ws.Range("K73").Value = 6352.799999999999
$ws.Range("L73").Value = 16791.201
$ws.Range("M73").Value = -5416.799999999999
$ws.Range("N73").Value = -18663.201

# Row 80
$ws.Range("H80").Value = 1843.4445
$ws.Range("I80").Value = 523.6667
$ws.Range("J80").Value = 2503.3333
$ws.Range("K80").Value = 1571.0001
$ws.Range("L80").Value = 7509.999899999999
$ws.Range("M80").Value = -573.0001
$ws.Range("N80").Value = -9505.999899999999

# Row 83
$ws.Range("H83").Value = 1843.4445
$ws.Range("I83").Value = 523.6667
$ws.Range("J83").Value = 2503.3333
$ws.Range("K83").Value = 4713.0003
$ws.Range("L83").Value = 22529.9997
$ws.Range("M83").Value = 278.9997000000003
$ws.Range("N83").Value = -32513.9997

# Row 135
$ws.Range("H135").Value = 1226.3334
$ws.Range("I135").Value = 1340.1666
$ws.Range("J135").Value = 1112.5
$ws.Range("K135").Value = 12061.4994
$ws.Range("L135").Value = 10012.5
$ws.Range("M135").Value = -9526.499400000001
$ws.Range("N135").Value = -15082.5

# ----- Worksheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 141
$ws.Range("H141").Value = 74548
$ws.Range("J141").Value = 74548
$ws.Range("L141").Value = 74548
$ws.Range("N141").Value = -84908

# ----- Worksheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 94
$ws.Range("H94").Value = 631.6667
$ws.Range("I94").Value = 685
$ws.Range("J94").Value = 525
$ws.Range("K94").Value = 685
$ws.Range("L94").Value = 525
$ws.Range("M94").Value = -234
$ws.Range("N94").Value = -1427

# ----- Worksheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Range("H58").Value = 1955.091
$ws.Range("I58").Value = 2124
$ws.Range("J58").Value = 1504.6666
$ws.Range("K58").Value = 2124
$ws.Range("L58").Value = 1504.6666
$ws.Range("M58").Value = -1921
$ws.Range("N58").Value = -1910.6666

# Row 136
$ws.Range("H136").Value = 1955.091
$ws.Range("I136").Value = 2124
$ws.Range("J136").Value = 1504.6666
$ws.Range("K136").Value = 6372
$ws.Range("L136").Value = 4513.9998
$ws.Range("M136").Value = -3822
$ws.Range("N136").Value = -9613.9998

# Row 140
$ws.Range("H140").Value = 61269.547
$ws.Range("J140").Value = 61269.547
$ws.Range("L140").Value = 61269.547
$ws.Range("N140").Value = -71629.54699999999

# ----- Worksheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 8
$ws.Range("H8").Value = 98.333336
$ws.Range("I8").Value = 98.333336
$ws.Range("K8").Value = 295.000008
$ws.Range("M8").Value = -156.000008

# Row 46
$ws.Range("H46").Value = 967.3333
$ws.Range("J46").Value = 1400
$ws.Range("L46").Value = 4200
$ws.Range("N46").Value = -4382

# Row 64
$ws.Range("H64").Value = 1212.125
$ws.Range("I64").Value = 212
$ws.Range("J64").Value = 1278.8
$ws.Range("K64").Value = 636
$ws.Range("L64").Value = 3836.4
$ws.Range("M64").Value = -366
$ws.Range("N64").Value = -4376.4

# Row 67
$ws.Range("H67").Value = 1212.125
$ws.Range("I67").Value = 212
$ws.Range("J67").Value = 1278.8
$ws.Range("K67").Value = 636
$ws.Range("L67").Value = 3836.4
$ws.Range("M67").Value = 300
$ws.Range("N67").Value = -5708.4

# Row 107
$ws.Range("H107").Value = 760.41174
$ws.Range("J107").Value = 673.3333
$ws.Range("L107").Value = 2019.9999
$ws.Range("N107").Value = -5859.9999

# Row 131
$ws.Range("H131").Value = 1029
$ws.Range("I131").Value = 498.42856
$ws.Range("J131").Value = 1068.9354
$ws.Range("K131").Value = 1495.28568
$ws.Range("L131").Value = 3206.8062
$ws.Range("M131").Value = 3544.71432
$ws.Range("N131").Value = -13286.8062

# ----- Worksheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 132
$ws.Range("H132").Value = 5318.8438
$ws.Range("I132").Value = 6684.1333
$ws.Range("J132").Value = 4114.1763
$ws.Range("K132").Value = 20052.3999
$ws.Range("L132").Value = 12342.5289
$ws.Range("M132").Value = -17522.3999
$ws.Range("N132").Value = -17402.5289

# Row 136
$ws.Range("H136").Value = 9261182
$ws.Range("I136").Value = 934.65
$ws.Range("J136").Value = 20836492
$ws.Range("K136").Value = 2803.95
$ws.Range("L136").Value = 62509476
$ws.Range("M136").Value = -253.9499999999998
$ws.Range("N136").Value = -62514576

# ----- Worksheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("N65").ClearContents()

# Row 81
$ws.Range("H81").Value = 766.5
$ws.Range("I81").Value = 636.7778
$ws.Range("K81").Value = 1273.5556
$ws.Range("M81").Value = -212.5555999999999

# Row 84
$ws.Range("H84").Value = 766.5
$ws.Range("I84").Value = 636.7778
$ws.Range("K84").Value = 6367.777999999999
$ws.Range("M84").Value = -1063.777999999999

# Row 136
$ws.Range("H136").Value = 3790.4546
$ws.Range("I136").Value = 610.4
$ws.Range("J136").Value = 7974.737
$ws.Range("K136").Value = 1831.2
$ws.Range("L136").Value = 23924.211
$ws.Range("M136").Value = 718.8000000000002
$ws.Range("N136").Value = -29024.211
